# Update cryptocurrency price/volume data to the latest scrape.
# Values are written with a leading apostrophe to force text storage
# (prices/percentages are display strings, not numeric types), then
# the style is reset to "Normal" so no extra number-format style is
# introduced on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.432.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.66%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.346.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -7.78%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'185.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -8.25%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'525.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -9.22%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -4.49%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.334.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -8.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.06%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.620"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -9.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'57.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.04%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -11.98%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -11.59%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -11.13%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.877.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "'  -4.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.344.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'64.117.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.93%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'17.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -11.82%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -12.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -11.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'371.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -9.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'PancakeSwap"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'3.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -13.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'80.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.57%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'10.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -18.94%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'5.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.84%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.86%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -10.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -11.53%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'8.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -11.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'651.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'28.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.81%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -14.85%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'11.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -10.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'59.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -7.10%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -10.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'35.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -14.37%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.373"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.32%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -9.44%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -15.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.767.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -13.37%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'PEPE"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0₃0614"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -20.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'WEMIXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -10.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0383"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -8.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -15.63%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -6.35%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'134.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.32%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.77%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -10.10%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Updated cryptos list"
